# Update "想去人数" (column F) values on both the "展览" and "全部类型"
# worksheets to reflect the newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    2  = 308
    3  = 13861
    7  = 278
    8  = 497
    10 = 89
    14 = 451
    15 = 5819
    16 = 133
    17 = 90
    18 = 979
    19 = 102
    20 = 56
    21 = 154
    22 = 262
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
